$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: split the run containing [$splitAt, $splitAt) by inserting a literal
# "-" there, then nudge the new run's font color away-and-back so the engine
# does not silently re-merge it with its (formatting-identical) neighbours.
# ---------------------------------------------------------------------------
function Split-WithDash($splitAt) {
    $rMid = $d.Range($splitAt, $splitAt)
    $rMid.InsertBefore("-")
    $rDash = $d.Range($splitAt, $splitAt + 1)
    $savedColor = $rDash.Font.Color
    $rDash.Font.Color = $savedColor + 1
    $rDash.Font.Color = $savedColor
}

# 1) "{% if languages %}"  ->  "{%" / "-" / " if languages %}"
$rng = $d.Content
$rng.Find.Execute("{% if languages %}")
$start = $rng.Start
Split-WithDash ($start + 2)

# 2) "{% if skills %}"  ->  "{%" / "-" / " if skills %}"
$rng = $d.Content
$rng.Find.Execute("{% if skills %}")
$start = $rng.Start
Split-WithDash ($start + 2)

# 3) " if education %}"  ->  " if education " / "-" / "%}"
$rng = $d.Content
$rng.Find.Execute(" if education %}")
$start = $rng.Start
$end = $rng.End
Split-WithDash ($end - 2)

# 4) " if hobbies %}"  ->  " if hobbies " / "-" / "%}"
$rng = $d.Content
$rng.Find.Execute(" if hobbies %}")
$start = $rng.Start
$end = $rng.End
Split-WithDash ($end - 2)

Write-Output "done"
